$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 17:22"

$ws.Range("A4").Value = "Estados Unidos"; $ws.Range("B4").Value = 616168; $ws.Range("C4").Value = 2282; $ws.Range("D4").Value = 40086; $ws.Range("E4").Value = 549882; $ws.Range("F4").Value = 13473; $ws.Range("G4").Value = 153; $ws.Range("H4").Value = 26200
$ws.Range("A8").Value = "Alemania"; $ws.Range("B8").Value = 132592; $ws.Range("C8").Value = 382; $ws.Range("D8").Value = 72600; $ws.Range("E8").Value = 56400; $ws.Range("F8").Value = 4288; $ws.Range("G8").Value = 97; $ws.Range("H8").Value = 3592
$ws.Range("A14").Value = "Paises Bajos"; $ws.Range("B14").Value = 28153; $ws.Range("C14").Value = 734; $ws.Range("D14").Value = 250; $ws.Range("E14").Value = 24769; $ws.Range("F14").Value = 1279; $ws.Range("G14").Value = 189; $ws.Range("H14").Value = 3134
$ws.Range("A16").Value = "Suiza"; $ws.Range("B16").Value = 26336; $ws.Range("C16").Value = 400; $ws.Range("D16").Value = 14700; $ws.Range("E16").Value = 10410; $ws.Range("F16").Value = 386; $ws.Range("G16").Value = 52; $ws.Range("H16").Value = 1226
$ws.Range("A23").Value = "India"; $ws.Range("B23").Value = 11555; $ws.Range("C23").Value = 68; $ws.Range("D23").Value = 1432; $ws.Range("E23").Value = 9718; $ws.Range("F23").Value = 0; $ws.Range("G23").Value = 12; $ws.Range("H23").Value = 405
$ws.Range("A27").Value = "Chile"; $ws.Range("B27").Value = 8273; $ws.Range("C27").Value = 356; $ws.Range("D27").Value = 2937; $ws.Range("E27").Value = 5242; $ws.Range("F27").Value = 389; $ws.Range("G27").Value = 2; $ws.Range("H27").Value = 94
$ws.Range("A28").Value = "Japon"; $ws.Range("B28").Value = 8100; $ws.Range("C28").Value = 215; $ws.Range("D28").Value = 853; $ws.Range("E28").Value = 7101; $ws.Range("F28").Value = 152; $ws.Range("G28").Value = 0; $ws.Range("H28").Value = 146
$ws.Range("A47").Value = "Singapur"; $ws.Range("B47").Value = 3699; $ws.Range("C47").Value = 447; $ws.Range("D47").Value = 652; $ws.Range("E47").Value = 3037; $ws.Range("F47").Value = 29; $ws.Range("G47").Value = 0; $ws.Range("H47").Value = 10
$ws.Range("A48").Value = "Republica Dominicana"; $ws.Range("B48").Value = 3614; $ws.Range("C48").Value = 328; $ws.Range("D48").Value = 208; $ws.Range("E48").Value = 3217; $ws.Range("F48").Value = 143; $ws.Range("G48").Value = 6; $ws.Range("H48").Value = 189
$ws.Range("A49").Value = "Panama"; $ws.Range("B49").Value = 3574; $ws.Range("C49").Value = 0; $ws.Range("D49").Value = 72; $ws.Range("E49").Value = 3407; $ws.Range("F49").Value = 106; $ws.Range("G49").Value = 0; $ws.Range("H49").Value = 95
$ws.Range("A50").Value = "Luxemburgo"; $ws.Range("B50").Value = 3307; $ws.Range("C50").Value = 0; $ws.Range("D50").Value = 500; $ws.Range("E50").Value = 2740; $ws.Range("F50").Value = 30; $ws.Range("G50").Value = 0; $ws.Range("H50").Value = 67
$ws.Range("A57").Value = "Grecia"; $ws.Range("B57").Value = 2192; $ws.Range("C57").Value = 22; $ws.Range("D57").Value = 269; $ws.Range("E57").Value = 1821; $ws.Range("F57").Value = 72; $ws.Range("G57").Value = 1; $ws.Range("H57").Value = 102
$ws.Range("A82").Value = "Cuba"; $ws.Range("B82").Value = 814; $ws.Range("C82").Value = 48; $ws.Range("D82").Value = 151; $ws.Range("E82").Value = 639; $ws.Range("F82").Value = 9; $ws.Range("G82").Value = 3; $ws.Range("H82").Value = 24
$ws.Range("A83").Value = "Afganistan"; $ws.Range("B83").Value = 784; $ws.Range("C83").Value = 70; $ws.Range("D83").Value = 43; $ws.Range("E83").Value = 716; $ws.Range("F83").Value = 0; $ws.Range("G83").Value = 2; $ws.Range("H83").Value = 25
$ws.Range("A85").Value = "Bulgaria"; $ws.Range("B85").Value = 747; $ws.Range("C85").Value = 34; $ws.Range("D85").Value = 105; $ws.Range("E85").Value = 606; $ws.Range("F85").Value = 31; $ws.Range("G85").Value = 1; $ws.Range("H85").Value = 36
$ws.Range("A99").Value = "Republica de Yibuti"; $ws.Range("B99").Value = 435; $ws.Range("C99").Value = 72; $ws.Range("D99").Value = 71; $ws.Range("E99").Value = 362; $ws.Range("F99").Value = 0; $ws.Range("G99").Value = 0; $ws.Range("H99").Value = 2
$ws.Range("A100").Value = "Honduras"; $ws.Range("B100").Value = 419; $ws.Range("C100").Value = 12; $ws.Range("D100").Value = 9; $ws.Range("E100").Value = 379; $ws.Range("F100").Value = 10; $ws.Range("G100").Value = 5; $ws.Range("H100").Value = 31
$ws.Range("A101").Value = "Guinea"; $ws.Range("B101").Value = 404; $ws.Range("C101").Value = 41; $ws.Range("D101").Value = 31; $ws.Range("E101").Value = 372; $ws.Range("F101").Value = 0; $ws.Range("G101").Value = 1; $ws.Range("H101").Value = 1
$ws.Range("A102").Value = "Malta"; $ws.Range("B102").Value = 399; $ws.Range("C102").Value = 6; $ws.Range("D102").Value = 44; $ws.Range("E102").Value = 352; $ws.Range("F102").Value = 4; $ws.Range("G102").Value = 0; $ws.Range("H102").Value = 3
$ws.Range("A103").Value = "Bolivia"; $ws.Range("B103").Value = 397; $ws.Range("C103").Value = 43; $ws.Range("D103").Value = 7; $ws.Range("E103").Value = 362; $ws.Range("F103").Value = 3; $ws.Range("G103").Value = 0; $ws.Range("H103").Value = 28
$ws.Range("A104").Value = "Jordania"; $ws.Range("B104").Value = 397; $ws.Range("C104").Value = 0; $ws.Range("D104").Value = 235; $ws.Range("E104").Value = 155; $ws.Range("F104").Value = 5; $ws.Range("G104").Value = 0; $ws.Range("H104").Value = 7
$ws.Range("A105").Value = "Taiwan"; $ws.Range("B105").Value = 395; $ws.Range("C105").Value = 2; $ws.Range("D105").Value = 137; $ws.Range("E105").Value = 252; $ws.Range("F105").Value = 0; $ws.Range("G105").Value = 0; $ws.Range("H105").Value = 6
$ws.Range("A106").Value = "Reunion"; $ws.Range("B106").Value = 391; $ws.Range("C106").Value = 0; $ws.Range("D106").Value = 40; $ws.Range("E106").Value = 351; $ws.Range("F106").Value = 3; $ws.Range("G106").Value = 0; $ws.Range("H106").Value = 0
$ws.Range("A107").Value = "Nigeria"; $ws.Range("B107").Value = 373; $ws.Range("C107").Value = 0; $ws.Range("D107").Value = 99; $ws.Range("E107").Value = 263; $ws.Range("F107").Value = 2; $ws.Range("G107").Value = 0; $ws.Range("H107").Value = 11
$ws.Range("A108").Value = "San Marino"; $ws.Range("B108").Value = 372; $ws.Range("C108").Value = 0; $ws.Range("D108").Value = 53; $ws.Range("E108").Value = 283; $ws.Range("F108").Value = 15; $ws.Range("G108").Value = 0; $ws.Range("H108").Value = 36
$ws.Range("A113").Value = "Montenegro"; $ws.Range("B113").Value = 288; $ws.Range("C113").Value = 5; $ws.Range("D113").Value = 55; $ws.Range("E113").Value = 229; $ws.Range("F113").Value = 7; $ws.Range("G113").Value = 0; $ws.Range("H113").Value = 4
$ws.Range("A132").Value = "Congo"; $ws.Range("B132").Value = 117; $ws.Range("C132").Value = 43; $ws.Range("D132").Value = 10; $ws.Range("E132").Value = 102; $ws.Range("F132").Value = 0; $ws.Range("G132").Value = 0; $ws.Range("H132").Value = 5
$ws.Range("A133").Value = "Trinidad yTobago"; $ws.Range("B133").Value = 114; $ws.Range("C133").Value = 1; $ws.Range("D133").Value = 19; $ws.Range("E133").Value = 87; $ws.Range("F133").Value = 0; $ws.Range("G133").Value = 0; $ws.Range("H133").Value = 8
$ws.Range("A134").Value = "Madagascar"; $ws.Range("B134").Value = 110; $ws.Range("C134").Value = 2; $ws.Range("D134").Value = 29; $ws.Range("E134").Value = 81; $ws.Range("F134").Value = 1; $ws.Range("G134").Value = 0; $ws.Range("H134").Value = 0
$ws.Range("A135").Value = "Jamaica"; $ws.Range("B135").Value = 105; $ws.Range("C135").Value = 0; $ws.Range("D135").Value = 21; $ws.Range("E135").Value = 79; $ws.Range("F135").Value = 0; $ws.Range("G135").Value = 1; $ws.Range("H135").Value = 5
$ws.Range("A136").Value = "Monaco"; $ws.Range("B136").Value = 93; $ws.Range("C136").Value = 0; $ws.Range("D136").Value = 6; $ws.Range("E136").Value = 86; $ws.Range("F136").Value = 5; $ws.Range("G136").Value = 0; $ws.Range("H136").Value = 1
$ws.Range("A137").Value = "Aruba"; $ws.Range("B137").Value = 92; $ws.Range("C137").Value = 0; $ws.Range("D137").Value = 32; $ws.Range("E137").Value = 59; $ws.Range("F137").Value = 1; $ws.Range("G137").Value = 1; $ws.Range("H137").Value = 1
$ws.Range("A138").Value = "Tanzania"; $ws.Range("B138").Value = 88; $ws.Range("C138").Value = 35; $ws.Range("D138").Value = 7; $ws.Range("E138").Value = 78; $ws.Range("F138").Value = 0; $ws.Range("G138").Value = 0; $ws.Range("H138").Value = 3
$ws.Range("A139").Value = "Guayana Francesa"; $ws.Range("B139").Value = 86; $ws.Range("C139").Value = 0; $ws.Range("D139").Value = 51; $ws.Range("E139").Value = 35; $ws.Range("F139").Value = 1; $ws.Range("G139").Value = 0; $ws.Range("H139").Value = 0
$ws.Range("A140").Value = "Etiopia"; $ws.Range("B140").Value = 85; $ws.Range("C140").Value = 3; $ws.Range("D140").Value = 15; $ws.Range("E140").Value = 67; $ws.Range("F140").Value = 0; $ws.Range("G140").Value = 0; $ws.Range("H140").Value = 3
$ws.Range("A141").Value = "Togo"; $ws.Range("B141").Value = 81; $ws.Range("C141").Value = 4; $ws.Range("D141").Value = 35; $ws.Range("E141").Value = 43; $ws.Range("F141").Value = 0; $ws.Range("G141").Value = 0; $ws.Range("H141").Value = 3
$ws.Range("A142").Value = "Gabon"; $ws.Range("B142").Value = 80; $ws.Range("C142").Value = 23; $ws.Range("D142").Value = 4; $ws.Range("E142").Value = 75; $ws.Range("F142").Value = 0; $ws.Range("G142").Value = 0; $ws.Range("H142").Value = 1
$ws.Range("A143").Value = "Liechtenstein"; $ws.Range("B143").Value = 79; $ws.Range("C143").Value = 0; $ws.Range("D143").Value = 55; $ws.Range("E143").Value = 23; $ws.Range("F143").Value = 0; $ws.Range("G143").Value = 0; $ws.Range("H143").Value = 1
$ws.Range("A144").Value = "Birmania"; $ws.Range("B144").Value = 74; $ws.Range("C144").Value = 11; $ws.Range("D144").Value = 2; $ws.Range("E144").Value = 68; $ws.Range("F144").Value = 0; $ws.Range("G144").Value = 0; $ws.Range("H144").Value = 4
$ws.Range("A145").Value = "Barbados"; $ws.Range("B145").Value = 73; $ws.Range("C145").Value = 0; $ws.Range("D145").Value = 15; $ws.Range("E145").Value = 53; $ws.Range("F145").Value = 4; $ws.Range("G145").Value = 0; $ws.Range("H145").Value = 5
$ws.Range("A146").Value = "Somalia"; $ws.Range("B146").Value = 60; $ws.Range("C146").Value = 0; $ws.Range("D146").Value = 2; $ws.Range("E146").Value = 56; $ws.Range("F146").Value = 2; $ws.Range("G146").Value = 0; $ws.Range("H146").Value = 2
$ws.Range("A169").Value = "Mozambique"; $ws.Range("B169").Value = 29; $ws.Range("C169").Value = 1; $ws.Range("D169").Value = 2; $ws.Range("E169").Value = 27; $ws.Range("F169").Value = 0; $ws.Range("G169").Value = 0; $ws.Range("H169").Value = 0
$ws.Range("A170").Value = "Siria"; $ws.Range("B170").Value = 29; $ws.Range("C170").Value = 0; $ws.Range("D170").Value = 5; $ws.Range("E170").Value = 22; $ws.Range("F170").Value = 0; $ws.Range("G170").Value = 0; $ws.Range("H170").Value = 2
$ws.Range("A194").Value = "Islas Malvinas"; $ws.Range("B194").Value = 11; $ws.Range("C194").Value = 0; $ws.Range("D194").Value = 1; $ws.Range("E194").Value = 10; $ws.Range("F194").Value = 0; $ws.Range("G194").Value = 0; $ws.Range("H194").Value = 0
$ws.Range("A195").Value = "Montserrat"; $ws.Range("B195").Value = 11; $ws.Range("C195").Value = 0; $ws.Range("D195").Value = 1; $ws.Range("E195").Value = 10; $ws.Range("F195").Value = 1; $ws.Range("G195").Value = 0; $ws.Range("H195").Value = 0